# Update "想去人数" (attendance/interest count) figures in the "展览" and
# "全部类型" worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7381
$ws1.Range("F5").Value = 274
$ws1.Range("F6").Value = 441
$ws1.Range("F7").Value = 3900
$ws1.Range("F8").Value = 313
$ws1.Range("F10").Value = 270
$ws1.Range("F11").Value = 620
$ws1.Range("F12").Value = 116

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7381
$ws4.Range("F7").Value = 274
$ws4.Range("F8").Value = 441
$ws4.Range("F9").Value = 3900
$ws4.Range("F10").Value = 313
$ws4.Range("F12").Value = 270
$ws4.Range("F13").Value = 620
$ws4.Range("F14").Value = 116
